# Applies the "library_layout" field-template addition to the
# slideseq-metadata workbook (HuBMAP docs/slideseq/slideseq-metadata.xlsx):
#
#   1. The existing "library_final_yield_unit list" lookup sheet is renamed
#      to "library_layout list" (so it keeps sheetId 6) and re-populated
#      with the single-end / paired-end options.
#   2. A brand-new sheet is inserted right after it, named
#      "library_final_yield_unit list" (becomes sheetId 7), re-populated
#      with the original "ng" option.
#   3. The comment on the library_layout header cell (S1) is reworded.
#   4. A list data-validation is added for the library_layout column (S),
#      positioned (like upstream) right after the N-column boolean
#      validation and before the V-column one.

$wb  = $excel.ActiveWorkbook
$tsv = $wb.Worksheets.Item("Export as TSV")

# --- 1 & 2: shuffle / create the lookup-list sheets ------------------------
$layoutSheet = $wb.Worksheets.Item("library_final_yield_unit list")
$layoutSheet.Name = "library_layout list"

$yieldUnitSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $layoutSheet)
$yieldUnitSheet.Name = "library_final_yield_unit list"

# Restore the yield-unit list's original content on the new sheet.
$yieldUnitSheet.Range("A1").Value = "ng"

# Populate the new library_layout list sheet.
$layoutSheet.Range("A1").Value = "single-end"
$layoutSheet.Range("A2").Value = "paired-end"

# --- 3: update the comment text on S1 --------------------------------------
$s1 = $tsv.Range("S1")
[void]$s1.Comment.Text("State whether the library was generated for single-end or paired end sequencing.")

# --- 4: add the data validation for the library_layout column (S) ----------
# The columns after N (V, AB, AC, AG, AH) already carry validations; remove
# them and re-add them in order, with the new S validation slotted in right
# after N (matching upstream's column ordering in the dataValidations list).
$vRange  = $tsv.Range("V2:V1048576")
$abRange = $tsv.Range("AB2:AB1048576")
$acRange = $tsv.Range("AC2:AC1048576")
$agRange = $tsv.Range("AG2:AG1048576")
$ahRange = $tsv.Range("AH2:AH1048576")

$vRange.Validation.Delete()
$abRange.Validation.Delete()
$acRange.Validation.Delete()
$agRange.Validation.Delete()
$ahRange.Validation.Delete()

$sRange = $tsv.Range("S2:S1048576")
$sValidation = $sRange.Validation
$sValidation.Add(3, 1, [System.Reflection.Missing]::Value, "'library_layout list'!`$A`$1:`$A`$2")
$sValidation.ErrorTitle = "Value must come from list"
$sValidation.ErrorMessage = "Value must be one of: single-end / paired-end."
$sValidation.ShowInput = $true
$sValidation.ShowError = $true

$vValidation = $vRange.Validation
$vValidation.Add(3, 1, [System.Reflection.Missing]::Value, '"TRUE,FALSE"')
$vValidation.ErrorTitle = "Not a boolean"
$vValidation.ErrorMessage = 'The values in this column must be "TRUE" or "FALSE".'
$vValidation.ShowInput = $true
$vValidation.ShowError = $true

$abValidation = $abRange.Validation
$abValidation.Add(2, 1, [System.Reflection.Missing]::Value, "-1e+307", "1e+307")
$abValidation.ErrorTitle = "Not a number"
$abValidation.ErrorMessage = "The values in this column must be numbers."
$abValidation.ShowInput = $true
$abValidation.ShowError = $true

$acValidation = $acRange.Validation
$acValidation.Add(3, 1, [System.Reflection.Missing]::Value, "'library_final_yield_unit list'!`$A`$1:`$A`$1")
$acValidation.ErrorTitle = "Value must come from list"
$acValidation.ErrorMessage = "Value must be one of: ng."
$acValidation.ShowInput = $true
$acValidation.ShowError = $true

$agValidation = $agRange.Validation
$agValidation.Add(2, 1, [System.Reflection.Missing]::Value, "-1e+307", "1e+307")
$agValidation.ErrorTitle = "Not a number"
$agValidation.ErrorMessage = "The values in this column must be numbers."
$agValidation.ShowInput = $true
$agValidation.ShowError = $true

$ahValidation = $ahRange.Validation
$ahValidation.Add(2, 1, [System.Reflection.Missing]::Value, "-1e+307", "1e+307")
$ahValidation.ErrorTitle = "Not a number"
$ahValidation.ErrorMessage = "The values in this column must be numbers."
$ahValidation.ShowInput = $true
$ahValidation.ShowError = $true

Write-Output "done"
